# Generate Report for Handback
# Update the "Correspond Handoff Datetime" / "Correspond Handback DateTime"
# values for the 253d68e7-... handback row on both the zh-cn and de-de
# report sheets to reflect the regenerated report timestamps.

$wb = $excel.ActiveWorkbook

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("E3").Value = "2016-03-19 07:51:12"
$zhcn.Range("H3").Value = "2016-03-19 07:51:54"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("E3").Value = "2016-03-19 07:51:20"
$dede.Range("H3").Value = "2016-03-19 07:52:07"
